$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 25.23990433333333
$ws.Range("H2").Value = 75.719713
$ws.Range("I2").Value = 0.05173702626903214
$ws.Range("J2").Value = 0.05173702626903214
$ws.Range("M2").Value = 0.110552
$ws.Range("N2").Value = 0.331656
$ws.Range("O2").Value = 0.01126249561724847
$ws.Range("P2").Value = 0.01126249561724847
$ws.Range("Q2").Value = 2.790321903858667
$ws.Range("R2").Value = 25.112897134728
$ws.Range("S2").Value = 0.0005826880316044435
$ws.Range("T2").Value = 0.0005826880316044436
$ws.Range("G3").Value = 25.23990433333333
$ws.Range("H3").Value = 75.719713
$ws.Range("I3").Value = 0.05173702626903214
$ws.Range("J3").Value = 0.05173702626903214
$ws.Range("O3").Value = 0.9181055646724333
$ws.Range("P3").Value = 0.9181055646724334
$ws.Range("Q3").Value = 227.4638014719062
$ws.Range("R3").Value = 2047.174213247156
$ws.Range("S3").Value = 0.04750005171720226
$ws.Range("T3").Value = 0.04750005171720228
$ws.Range("G4").Value = 25.23990433333333
$ws.Range("H4").Value = 75.719713
$ws.Range("I4").Value = 0.05173702626903214
$ws.Range("J4").Value = 0.05173702626903214
$ws.Range("M4").Value = 0.6933189999999999
$ws.Range("N4").Value = 2.079957
$ws.Range("O4").Value = 0.07063193971031816
$ws.Range("P4").Value = 0.07063193971031817
$ws.Range("Q4").Value = 17.49930523248233
$ws.Range("R4").Value = 157.493747092341
$ws.Range("S4").Value = 0.003654286520225424
$ws.Range("T4").Value = 0.003654286520225426
$ws.Range("I5").Value = 0.8454897015965644
$ws.Range("J5").Value = 0.8454897015965646
$ws.Range("M5").Value = 0.110552
$ws.Range("N5").Value = 0.331656
$ws.Range("O5").Value = 0.01126249561724847
$ws.Range("P5").Value = 0.01126249561724847
$ws.Range("Q5").Value = 45.59961412517333
$ws.Range("R5").Value = 410.39652712656
$ws.Range("S5").Value = 0.009522324058660025
$ws.Range("T5").Value = 0.009522324058660029
$ws.Range("I6").Value = 0.8454897015965644
$ws.Range("J6").Value = 0.8454897015965646
$ws.Range("O6").Value = 0.9181055646724333
$ws.Range("P6").Value = 0.9181055646724334
$ws.Range("S6").Value = 0.7762487999090409
$ws.Range("T6").Value = 0.7762487999090412
$ws.Range("I7").Value = 0.8454897015965644
$ws.Range("J7").Value = 0.8454897015965646
$ws.Range("M7").Value = 0.6933189999999999
$ws.Range("N7").Value = 2.079957
$ws.Range("O7").Value = 0.07063193971031816
$ws.Range("P7").Value = 0.07063193971031817
$ws.Range("Q7").Value = 285.9747346556466
$ws.Range("R7").Value = 2573.77261190082
$ws.Range("S7").Value = 0.05971857762886343
$ws.Range("T7").Value = 0.05971857762886345
$ws.Range("G8").Value = 50.137933
$ws.Range("H8").Value = 150.413799
$ws.Range("I8").Value = 0.1027732721344034
$ws.Range("J8").Value = 0.1027732721344034
$ws.Range("M8").Value = 0.110552
$ws.Range("N8").Value = 0.331656
$ws.Range("O8").Value = 0.01126249561724847
$ws.Range("P8").Value = 0.01126249561724847
$ws.Range("Q8").Value = 5.542848769016
$ws.Range("R8").Value = 49.88563892114399
$ws.Range("S8").Value = 0.001157483526984002
$ws.Range("T8").Value = 0.001157483526984003
$ws.Range("G9").Value = 50.137933
$ws.Range("H9").Value = 150.413799
$ws.Range("I9").Value = 0.1027732721344034
$ws.Range("J9").Value = 0.1027732721344034
$ws.Range("O9").Value = 0.9181055646724333
$ws.Range("P9").Value = 0.9181055646724334
$ws.Range("Q9").Value = 451.8465952765986
$ws.Range("R9").Value = 4066.619357489387
$ws.Range("S9").Value = 0.09435671304619006
$ws.Range("T9").Value = 0.09435671304619006
$ws.Range("G10").Value = 50.137933
$ws.Range("H10").Value = 150.413799
$ws.Range("I10").Value = 0.1027732721344034
$ws.Range("J10").Value = 0.1027732721344034
$ws.Range("M10").Value = 0.6933189999999999
$ws.Range("N10").Value = 2.079957
$ws.Range("O10").Value = 0.07063193971031816
$ws.Range("P10").Value = 0.07063193971031817
$ws.Range("Q10").Value = 34.76158156962699
$ws.Range("R10").Value = 312.854234126643
$ws.Range("S10").Value = 0.007259075561229299
$ws.Range("T10").Value = 0.0072590755612293
